$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8641925454139709
$ws.Range("B1").Value = 2.611090183258057
$ws.Range("C1").Value = 1.262760877609253
$ws.Range("D1").Value = 1.265014886856079
$ws.Range("E1").Value = 1.394752740859985
